$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Scans" / table "Scans": append a new scan row (2026-01-27 10:21:34)
# with no new publication IDs found.
# ---------------------------------------------------------------------------
$wsScans = $wb.Worksheets.Item("Scans")
$tblScans = $wsScans.ListObjects.Item("Scans")
$tblScans.ListRows.Add() | Out-Null

# "2026-01-27" already exists elsewhere in this column as text (shared string).
# Assigning it directly would be auto-recognised as a date, so copy the
# existing text cell instead to preserve its text type.
$wsScans.Range("A2").Copy() | Out-Null
$wsScans.Range("A4").PasteSpecial() | Out-Null
$wsScans.Cells.Item(4, 2).Value = "10:21:34"

# ---------------------------------------------------------------------------
# Sheet "Order Papers" / table "Order_Papers": populate the "HC matched"
# column for the existing rows, then append a duplicated block of rows
# (same four distinct reports) reflecting the new order paper scan.
# ---------------------------------------------------------------------------
$wsOP = $wb.Worksheets.Item("Order Papers")
$tblOP = $wsOP.ListObjects.Item("Order_Papers")

$wsOP.Cells.Item(2, 7).Value = "Published"
$wsOP.Cells.Item(3, 7).Value = "Missing"
$wsOP.Cells.Item(4, 7).Value = "Missing"
$wsOP.Cells.Item(5, 7).Value = "Missing"
$wsOP.Cells.Item(6, 7).Value = "Published"
$wsOP.Cells.Item(7, 7).Value = "Missing"
$wsOP.Cells.Item(8, 7).Value = "Missing"
$wsOP.Cells.Item(9, 7).Value = "Missing"

for ($i = 0; $i -lt 4; $i++) {
  $tblOP.ListRows.Add() | Out-Null
}

$newRows = @(
  @(10, "Health and Social Care", "5th Report: First 1000 Days: a renewed focus", "HC 802", "00:01:00", "Published"),
  @(11, "International Development", "7th Special Report: Empowering Development: Energy Access for Communities: Government Response", "HC 1626", "00:01:00", "Missing"),
  @(12, "Treasury", "6th Special Report: Taxation of gambling: Government Response", "HC 1625", "00:01:00", "Missing"),
  @(13, "Home Affairs", "3rd Special Report: The Home Office`u{2019}s management of asylum accommodation: Government Response", "HC 1642", "10:00:00", "Missing")
)

foreach ($row in $newRows) {
  $r = $row[0]

  # Column A & E hold "2026-01-22" which already exists as text elsewhere in
  # the sheet; copy it across so it is not re-interpreted as a date serial.
  $wsOP.Range("A2").Copy() | Out-Null
  $wsOP.Range("A$r").PasteSpecial() | Out-Null
  $wsOP.Range("A2").Copy() | Out-Null
  $wsOP.Range("E$r").PasteSpecial() | Out-Null

  $wsOP.Cells.Item($r, 2).Value = $row[1]
  $wsOP.Cells.Item($r, 3).Value = $row[2]
  $wsOP.Cells.Item($r, 4).Value = $row[3]
  $wsOP.Cells.Item($r, 6).Value = $row[4]
  $wsOP.Cells.Item($r, 7).Value = $row[5]
}
